$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10:27 down to 11:28
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44414
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 100112026
$ws.Range("G10").Value = "Haba"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 18000
$ws.Range("L10").Value = 18000
$ws.Range("M10").Value = 18000
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 720
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
